$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 59-81 each shift down by one (row N now holds what used to be in row N-1),
# row 59 receives a brand-new data point, and a new row 82 is appended holding
# what used to be in row 81.

$ws.Range("D59").Value = 45006
$ws.Range("I59").Value = 'Primera'
$ws.Range("J59").Value = 80
$ws.Range("K59").Value = 17000
$ws.Range("L59").Value = 18000
$ws.Range("M59").Value = 17500
$ws.Range("P59").Value = 972

$ws.Range("D60").Value = 44232
$ws.Range("I60").Value = 'Especial'
$ws.Range("J60").Value = 50
$ws.Range("K60").Value = 22000
$ws.Range("L60").Value = 22000
$ws.Range("M60").Value = 22000
$ws.Range("P60").Value = 1222

$ws.Range("D61").Value = 44691
$ws.Range("I61").Value = 'Primera'
$ws.Range("J61").Value = 70
$ws.Range("K61").Value = 17000
$ws.Range("L61").Value = 17000
$ws.Range("M61").Value = 17000
$ws.Range("P61").Value = 944

$ws.Range("D62").Value = 44637
$ws.Range("I62").Value = 'Especial'
$ws.Range("J62").Value = 50
$ws.Range("K62").Value = 21000
$ws.Range("L62").Value = 21000
$ws.Range("M62").Value = 21000
$ws.Range("P62").Value = 1167

$ws.Range("D63").Value = 44350
$ws.Range("I63").Value = 'Primera'
$ws.Range("J63").Value = 20
$ws.Range("K63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("M63").Value = 20000
$ws.Range("P63").Value = 1111

$ws.Range("D64").Value = 44721
$ws.Range("I64").Value = 'Especial'
$ws.Range("J64").Value = 30
$ws.Range("K64").Value = 21000
$ws.Range("L64").Value = 21000
$ws.Range("M64").Value = 21000
$ws.Range("P64").Value = 1167

$ws.Range("D65").Value = 44764
$ws.Range("I65").Value = 'Primera'
$ws.Range("J65").Value = 50
$ws.Range("K65").Value = 18000
$ws.Range("L65").Value = 18000
$ws.Range("M65").Value = 18000
$ws.Range("P65").Value = 1000

$ws.Range("D66").Value = 44749
$ws.Range("I66").Value = 'Especial'
$ws.Range("J66").Value = 35
$ws.Range("K66").Value = 21000
$ws.Range("L66").Value = 21000
$ws.Range("M66").Value = 21000
$ws.Range("P66").Value = 1167

$ws.Range("D67").Value = 44974
$ws.Range("I67").Value = 'Primera'
$ws.Range("J67").Value = 60
$ws.Range("K67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("M67").Value = 20000
$ws.Range("P67").Value = 1111

$ws.Range("D68").Value = 44371
$ws.Range("I68").Value = 'Especial'
$ws.Range("J68").Value = 20
$ws.Range("K68").Value = 20000
$ws.Range("L68").Value = 20000
$ws.Range("M68").Value = 20000
$ws.Range("P68").Value = 1111

$ws.Range("D69").Value = 44979
$ws.Range("I69").Value = 'Primera'
$ws.Range("J69").Value = 25
$ws.Range("K69").Value = 20000
$ws.Range("L69").Value = 20000
$ws.Range("M69").Value = 20000
$ws.Range("P69").Value = 1111

$ws.Range("D70").Value = 44238
$ws.Range("I70").Value = 'Especial'
$ws.Range("J70").Value = 50
$ws.Range("K70").Value = 20000
$ws.Range("L70").Value = 20000
$ws.Range("M70").Value = 20000
$ws.Range("P70").Value = 1111

$ws.Range("D71").Value = 44663
$ws.Range("I71").Value = 'Primera'
$ws.Range("J71").Value = 90
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 18000
$ws.Range("M71").Value = 18000
$ws.Range("P71").Value = 1000

$ws.Range("D72").Value = 44285
$ws.Range("I72").Value = 'Especial'
$ws.Range("J72").Value = 70
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 18000
$ws.Range("M72").Value = 18000
$ws.Range("P72").Value = 1000

$ws.Range("D73").Value = 44292
$ws.Range("I73").Value = 'Especial'
$ws.Range("J73").Value = 70
$ws.Range("K73").Value = 17000
$ws.Range("L73").Value = 17000
$ws.Range("M73").Value = 17000
$ws.Range("P73").Value = 944

$ws.Range("D74").Value = 44253
$ws.Range("I74").Value = 'Primera'
$ws.Range("J74").Value = 70
$ws.Range("K74").Value = 18000
$ws.Range("L74").Value = 18000
$ws.Range("M74").Value = 18000
$ws.Range("P74").Value = 1000

$ws.Range("D75").Value = 44278
$ws.Range("I75").Value = 'Especial'
$ws.Range("J75").Value = 70
$ws.Range("K75").Value = 18000
$ws.Range("L75").Value = 18000
$ws.Range("M75").Value = 18000
$ws.Range("P75").Value = 1000

$ws.Range("D76").Value = 44306
$ws.Range("I76").Value = 'Especial'
$ws.Range("J76").Value = 80
$ws.Range("K76").Value = 18000
$ws.Range("L76").Value = 18000
$ws.Range("M76").Value = 18000
$ws.Range("P76").Value = 1000

$ws.Range("D77").Value = 44257
$ws.Range("I77").Value = 'Primera'
$ws.Range("J77").Value = 60
$ws.Range("K77").Value = 16000
$ws.Range("L77").Value = 16000
$ws.Range("M77").Value = 16000
$ws.Range("P77").Value = 889

$ws.Range("D78").Value = 44364
$ws.Range("I78").Value = 'Especial'
$ws.Range("J78").Value = 30
$ws.Range("K78").Value = 20000
$ws.Range("L78").Value = 20000
$ws.Range("M78").Value = 20000
$ws.Range("P78").Value = 1111

$ws.Range("D79").Value = 44245
$ws.Range("I79").Value = 'Primera'
$ws.Range("J79").Value = 40
$ws.Range("K79").Value = 18000
$ws.Range("L79").Value = 18000
$ws.Range("M79").Value = 18000
$ws.Range("P79").Value = 1000

$ws.Range("D80").Value = 44242
$ws.Range("I80").Value = 'Especial'
$ws.Range("J80").Value = 50
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("M80").Value = 20000
$ws.Range("P80").Value = 1111

$ws.Range("D81").Value = 44988
$ws.Range("I81").Value = 'Primera'
$ws.Range("J81").Value = 100
$ws.Range("K81").Value = 20000
$ws.Range("L81").Value = 22000
$ws.Range("M81").Value = 21000
$ws.Range("P81").Value = 1167

$ws.Range("I82").Value = 'Especial'
$ws.Range("J82").Value = 60
$ws.Range("K82").Value = 20000
$ws.Range("L82").Value = 20000
$ws.Range("M82").Value = 20000
$ws.Range("P82").Value = 1111

# Row 82 did not exist before, so also populate the columns that are constant
# across every data row in this sheet, matching row 81 (and all others).
$ws.Range("A82").Value = 4
$ws.Range("B82").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C82").Value = 'Los Lagos'
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = 100112043
$ws.Range("G82").Value = 'Pepino dulce'
$ws.Range("H82").Value = 'Cultivar IV Región'
$ws.Range("N82").Value = '$/bandeja 18 kilos'
$ws.Range("O82").Value = 'Provincia de Limarí'
$ws.Range("Q82").Value = 18
$ws.Range("R82").Value = 'Hortaliza'

# D82 uses the same date number format as every other Fecha cell in the column.
$ws.Range("D82").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("D82").Value = 44236

